# Re-running the handoff-report generation for the 4 "Ready for handoff" rows
# (01eb3ba1…, 5ef03e35…, 66821751…, de8ccdf9…) bumps their priority from
# "low" to "ht", and refreshes the generated/handoff timestamps.
$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date is regenerated a little later.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-08-27 18:38:53"

# "zh-cn" sheet: Priority low -> ht, and the handoff xliff was regenerated,
# updating the Latest Handoff Datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-08-27 18:38:48"

# "de-de" sheet: same priority bump; its Latest Handoff Datetime tracks the
# same refreshed generation time shown on the Overview sheet.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-08-27 18:38:53"
